# Resolved css issue of data table on listing page.
# Insert a new attribute row ("valid_ip_addresses") into the
# "Software License Agreement" block, just above "terms and conditions",
# shifting all following rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 39 (old row 39 "terms and conditions" and
# everything below it shifts down to 40, 41, ... 51).
$ws.Range("A39").EntireRow.Insert() | Out-Null

# Fill in the new cell's text.
$ws.Range("C39").Value = "valid_ip_addresses"

# Give it its own (black) font so it gets a distinct style, matching the
# new font/cellXf entries added to styles.xml.
$ws.Range("C39").Font.Color = 0

# Restore the selection / scroll position to where the author left off.
$excel.ActiveWindow.ScrollRow = 28
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C39").Select() | Out-Null
